# Applies the recorded edit to the "artfynd" worksheet:
#  1. The data rows 2-9 (one record each) get reordered. Row 8 keeps its
#     place; the others are permuted as follows (new row <- old row):
#       2<-5, 3<-6, 4<-2, 5<-7, 6<-9, 7<-4, 8<-8, 9<-3
#  2. The "Ost"/"Nord" coordinate columns (Q, R) are rounded to whole
#     numbers for every data row.
#  3. The "Starttid"/"Sluttid" time columns (Z, AB) are cleared for every
#     data row (they only ever held the redundant "00:00" value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 9
$stageOffset = 1000

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")
$lastCol = $cols[$cols.Count - 1]

# Mapping of new row number -> original row number.
$rowMap = @{
    2 = 5
    3 = 6
    4 = 2
    5 = 7
    6 = 9
    7 = 4
    8 = 8
    9 = 3
}

# --- Step 1: remember exactly which cells are populated in each source
#     row before anything is touched, so the permutation below can
#     reproduce the same sparse layout (some rows have blank placeholder
#     cells such as J/N/AF while others simply omit them). ---
$populated = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $colsWithData = @()
    foreach ($col in $cols) {
        $v = $ws.Range($col + $r).Value2
        if ($v -ne $null) {
            $colsWithData += $col
        }
    }
    $populated[$r] = $colsWithData
}

# --- Step 2: stage a copy of every data row well out of the way so the
#     permutation below never overwrites data it still needs to read. ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $stageRow = $r + $stageOffset
    $dst = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
    $src.Copy($dst)
}

# --- Step 3: write each destination row from its staged source row, then
#     restore the original sparse cell layout that belongs with that
#     record (Copy() fills in every intervening column, so blank out the
#     ones that were not actually present originally). ---
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $stageRow = $oldRow + $stageOffset
    $src = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
    $dst = $ws.Range("A" + $newRow + ":" + $lastCol + $newRow)
    $src.Copy($dst)

    $keepCols = $populated[$oldRow]
    foreach ($col in $cols) {
        if (-not ($keepCols -contains $col)) {
            $ws.Range($col + $newRow).ClearContents()
        }
    }
}

# --- Step 4: clean up the staging area. ---
$stageFirst = $firstRow + $stageOffset
$stageLast = $lastRow + $stageOffset
$ws.Range("A" + $stageFirst + ":" + $lastCol + $stageLast).Clear()

# --- Step 5: round the Ost (Q) / Nord (R) coordinates for every row. ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $q = $ws.Range("Q" + $r).Value2
    if ($q -ne $null) {
        $ws.Range("Q" + $r).Value = [Math]::Round([double]$q)
    }
    $rr = $ws.Range("R" + $r).Value2
    if ($rr -ne $null) {
        $ws.Range("R" + $r).Value = [Math]::Round([double]$rr)
    }
}

# --- Step 6: drop the redundant Starttid (Z) / Sluttid (AB) values. ---
$ws.Range("Z" + $firstRow + ":Z" + $lastRow).ClearContents()
$ws.Range("AB" + $firstRow + ":AB" + $lastRow).ClearContents()
